$d = $word.ActiveDocument

# --- 1. Resize the "client/server" protocol table (table #3) ---
# tblW 5046 -> 7380 dxa ; third column (and its tcW cells) 1716 -> 4050 dxa.
# Word COM widths are expressed in points (1 pt = 20 dxa/twips).
$tbl = $d.Tables(3)
$tbl.PreferredWidth = 7380 / 20
$tbl.Columns(3).Width = 4050 / 20

# --- 2. Add the new 'location' line to the DONE message body ---
# Row 2 (the DONE message), column 3, paragraph 4 is "     'type' : DONE".
# Insert a new paragraph right after it (before the closing "}" paragraph).
$cell = $tbl.Cell(2, 3)
$donePara = $cell.Range.Paragraphs(4)
$insertPoint = $d.Range($donePara.Range.End, $donePara.Range.End)
$locationXml = '<w:p xmlns:w="http://schemas.openxmlformats.org/wordprocessingml/2006/main">' + `
  '<w:pPr><w:bidi/><w:jc w:val="right"/></w:pPr>' + `
  '<w:r><w:t xml:space="preserve">     </w:t></w:r>' + `
  '<w:r><w:t>‘</w:t></w:r>' + `
  '<w:r><w:t>location’</w:t></w:r>' + `
  '<w:r><w:t>: {</w:t></w:r>' + `
  '<w:r><w:t>‘row’ : &lt;num&gt;, ‘col’ : &lt;num&gt;}</w:t></w:r>' + `
  '</w:p>'
$insertPoint.InsertXML($locationXml) | Out-Null

# --- 3. Strip the stray RTL run-formatting from the empty paragraph after the last table ---
$lastTbl = $d.Tables($d.Tables.Count)
$afterTblRange = $d.Range($lastTbl.Range.End, $lastTbl.Range.End + 2)
$emptyPara = $afterTblRange.Paragraphs(1)
$emptyXml = '<w:p xmlns:w="http://schemas.openxmlformats.org/wordprocessingml/2006/main"/>'
$emptyPara.Range.InsertXML($emptyXml) | Out-Null
